$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.998.18"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "2.460.37"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'580.66"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").Value = "'166.35"
$ws.Range("E6").Value = "  -3.92%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").Value = "2.459.89"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("D13").Value = "'0.333"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").Value = "2.929.04"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "'25.35"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "66.831.80"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("E17").Value = "  -4.41%  "
$ws.Range("D18").Value = "2.442.50"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").Value = "'11.31"
$ws.Range("E19").Value = "  -3.78%  "
$ws.Range("D20").Value = "'7.65"
$ws.Range("E20").Value = "  -4.66%  "
$ws.Range("D21").Value = "'353.59"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").Value = "'4.01"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "'69.35"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").Value = "'4.22"
$ws.Range("E25").Value = "  -7.14%  "
$ws.Range("D26").Value = "'1.76"
$ws.Range("E26").Value = "  -7.34%  "
$ws.Range("D27").Value = "'8.94"
$ws.Range("E27").Value = "  -8.61%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "2.587.18"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "0.0₃0898"
$ws.Range("E30").Value = "  -6.63%  "
$ws.Range("D31").Value = "'507.48"
$ws.Range("E31").Value = "  -4.32%  "
$ws.Range("D32").Value = "'7.77"
$ws.Range("E32").Value = "  -5.39%  "
$ws.Range("D33").Value = "'1.78"
$ws.Range("E33").Value = "  -4.48%  "
$ws.Range("E34").Value = "  -5.67%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.117"
$ws.Range("E36").Value = "  -7.61%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'157.98"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'18.46"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'18.57"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  -5.93%  "
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  -6.23%  "
$ws.Range("D43").Value = "'0.327"
$ws.Range("E43").Value = "  -5.99%  "
$ws.Range("D44").Value = "'4.78"
$ws.Range("E44").Value = "  -6.37%  "
$ws.Range("D45").Value = "'38.68"
$ws.Range("D46").Value = "'2.32"
$ws.Range("E46").Value = "  -6.60%  "
$ws.Range("D47").Value = "'141.34"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").Value = "'3.47"
$ws.Range("E48").Value = "  -5.73%  "
$ws.Range("D49").Value = "'0.514"
$ws.Range("E49").Value = "  -5.83%  "
$ws.Range("E50").Value = "  -6.55%  "
$ws.Range("D51").Value = "'1.59"
$ws.Range("E51").Value = "  -5.57%  "
